# Matriz de trazabilidad - registrar 3 nuevos requerimientos (filas 8, 9, 10)
# relacionados con los cambios de tipos de datos (DATE -> TIME / DATETIME)
# mencionados en el commit ("Fixes y Avance de administrador").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Matriz de trazabilidad")
$ws.Activate()

# --- Fila 8 ---------------------------------------------------------------
$ws.Range("B8").Value = "horaInicial y horaFinal en la tabla cursos fueron cambiados de DATE a TIME"
$ws.Range("C8").Value = "Time es específico para horas"
$ws.Range("D6").Copy()
$ws.Range("D8").PasteSpecial(-4122)
$ws.Range("D8").Formula = "=TODAY()"
$ws.Range("E8").Value = "Cristian Prince"
$ws.Range("F8").Value = "Requisito"
$ws.Range("G8").Value = "Alta"
$ws.Range("H8").Value = "Implementado"
$ws.Range("I8").Value = "Requisito faltante"

# --- Fila 9 -----------------------------------------------------------------
$ws.Range("B9").Value = "fechaInicio y fechaFin en la tabla evaluaciones cambiaron de DATE a DATETIME"
$ws.Range("C9").Value = "DATETIME es más preciso para el día y hora"
$ws.Range("D6").Copy()
$ws.Range("D9").PasteSpecial(-4122)
$ws.Range("D9").Formula = "=TODAY()"
$ws.Range("E9").Value = "Cristian Prince"
$ws.Range("F9").Value = "Requisito"
$ws.Range("G9").Value = "Alta"
$ws.Range("H9").Value = "Implementado"
$ws.Range("I9").Value = "Requisito faltante"

# --- Fila 10 ------------------------------------------------------------
$ws.Range("B10").Value = "fechaEntrega en la tabla calificaciones fue cambiado de DATE a DATETIME"
$ws.Range("C10").Value = "DATETIME es más preciso para el día y hora"
$ws.Range("D6").Copy()
$ws.Range("D10").PasteSpecial(-4122)
$ws.Range("D10").Formula = "=TODAY()"
$ws.Range("E10").Value = "Cristian Prince"
$ws.Range("F10").Value = "Requisito"
$ws.Range("G10").Value = "Alta"
$ws.Range("H10").Value = "Implementado"
$ws.Range("I10").Value = "Requisito faltante"

# Last active selection left on the sheet by the author.
$ws.Range("F12").Select() | Out-Null
